$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before the existing "update_usr_id" column (H) to hold
# the new tenant_id field. Existing H/I (update_usr_id_lbl / update_time_lbl)
# shift right to I/J.
$ws.Columns("H").Insert()

$ws.Range("H1").Value = '<%=comment.tenant_id_lbl%><%selectList.tenant_id = data.findAllTenant.map((item) => item.lbl)%><%_dataValidation_({ sqref: `${ _col }2:${ _col }${ _lastRow }`, formula1: `"${ selectList.tenant_id.join(",") }"` })%>'
$ws.Range("H2").Value = '<%=model.tenant_id_lbl%>'
